$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete year rows (2000年, 2004年-2009年) that sat above the
# 2010年/2011年 rows so the surviving data shifts up to rows 2-3.
$ws.Rows("2:8").Delete()
